# Auto-generated edit script updating leve profit calculation cells
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook


$ws = $wb.Sheets.Item("ALC")
$ws.Cells.Item(10, 8).Value = 5362.75  # H10: 5483.6665 -> 5362.75
$ws.Cells.Item(10, 10).Value = 7000  # J10: 8000 -> 7000
$ws.Cells.Item(10, 12).Value = 7000  # L10: 8000 -> 7000
$ws.Cells.Item(10, 14).Value = -7586  # N10: -8586 -> -7586
$ws.Cells.Item(80, 8).Value = 5152.44  # H80: 5435.391 -> 5152.44
$ws.Cells.Item(80, 9).Value = 406.84616  # I80: 432.41666 -> 406.84616
$ws.Cells.Item(80, 10).Value = 10293.5  # J80: 10893.182 -> 10293.5
$ws.Cells.Item(80, 11).Value = 1220.53848  # K80: 1297.24998 -> 1220.53848
$ws.Cells.Item(80, 12).Value = 30880.5  # L80: 32679.546 -> 30880.5
$ws.Cells.Item(80, 13).Value = -222.5384799999999  # M80: -299.2499800000001 -> -222.5384799999999
$ws.Cells.Item(80, 14).Value = -32876.5  # N80: -34675.546 -> -32876.5
$ws.Cells.Item(83, 8).Value = 5152.44  # H83: 5435.391 -> 5152.44
$ws.Cells.Item(83, 9).Value = 406.84616  # I83: 432.41666 -> 406.84616
$ws.Cells.Item(83, 10).Value = 10293.5  # J83: 10893.182 -> 10293.5
$ws.Cells.Item(83, 11).Value = 3661.61544  # K83: 3891.74994 -> 3661.61544
$ws.Cells.Item(83, 12).Value = 92641.5  # L83: 98038.63800000001 -> 92641.5
$ws.Cells.Item(83, 13).Value = 1330.38456  # M83: 1100.25006 -> 1330.38456
$ws.Cells.Item(83, 14).Value = -102625.5  # N83: -108022.638 -> -102625.5
$ws.Cells.Item(86, 8).Value = 32152.758  # H86: 32201.152 -> 32152.758
$ws.Cells.Item(86, 9).Value = 73013.5  # I86: 73049 -> 73013.5
$ws.Cells.Item(86, 10).Value = 2044.8422  # J86: 2102.7368 -> 2044.8422
$ws.Cells.Item(86, 11).Value = 73013.5  # K86: 73049 -> 73013.5
$ws.Cells.Item(86, 12).Value = 2044.8422  # L86: 2102.7368 -> 2044.8422
$ws.Cells.Item(86, 13).Value = -71890.5  # M86: -71926 -> -71890.5
$ws.Cells.Item(86, 14).Value = -4290.8422  # N86: -4348.736800000001 -> -4290.8422
$ws.Cells.Item(89, 8).Value = 32152.758  # H89: 32201.152 -> 32152.758
$ws.Cells.Item(89, 9).Value = 73013.5  # I89: 73049 -> 73013.5
$ws.Cells.Item(89, 10).Value = 2044.8422  # J89: 2102.7368 -> 2044.8422
$ws.Cells.Item(89, 11).Value = 365067.5  # K89: 365245 -> 365067.5
$ws.Cells.Item(89, 12).Value = 10224.211  # L89: 10513.684 -> 10224.211
$ws.Cells.Item(89, 13).Value = -359451.5  # M89: -359629 -> -359451.5
$ws.Cells.Item(89, 14).Value = -21456.211  # N89: -21745.684 -> -21456.211
$ws.Cells.Item(98, 8).Value = 1629  # H98: 975.3333 -> 1629
$ws.Cells.Item(98, 9).Value = 1899.3636  # I98: 999.62964 -> 1899.3636
$ws.Cells.Item(98, 10).Value = 1133.3334  # J98: 866 -> 1133.3334
$ws.Cells.Item(98, 11).Value = 1899.3636  # K98: 999.62964 -> 1899.3636
$ws.Cells.Item(98, 12).Value = 1133.3334  # L98: 866 -> 1133.3334
$ws.Cells.Item(98, 13).Value = -401.3635999999999  # M98: 498.37036 -> -401.3635999999999
$ws.Cells.Item(98, 14).Value = -4129.3334  # N98: -3862 -> -4129.3334
$ws.Cells.Item(122, 8).Value = 1629  # H122: 975.3333 -> 1629
$ws.Cells.Item(122, 9).Value = 1899.3636  # I122: 999.62964 -> 1899.3636
$ws.Cells.Item(122, 10).Value = 1133.3334  # J122: 866 -> 1133.3334
$ws.Cells.Item(122, 11).Value = 5698.0908  # K122: 2998.88892 -> 5698.0908
$ws.Cells.Item(122, 12).Value = 3400.0002  # L122: 2598 -> 3400.0002
$ws.Cells.Item(122, 13).Value = -3248.0908  # M122: -548.8889199999999 -> -3248.0908
$ws.Cells.Item(122, 14).Value = -8300.0002  # N122: -7498 -> -8300.0002
$ws.Cells.Item(132, 8).Value = 1578.7391  # H132: 1979.9333 -> 1578.7391
$ws.Cells.Item(132, 9).Value = 1572.9062  # I132: 1668.6 -> 1572.9062
$ws.Cells.Item(132, 10).Value = 1592.0714  # J132: 2602.6 -> 1592.0714
$ws.Cells.Item(132, 11).Value = 4718.7186  # K132: 5005.799999999999 -> 4718.7186
$ws.Cells.Item(132, 12).Value = 4776.2142  # L132: 7807.799999999999 -> 4776.2142
$ws.Cells.Item(132, 13).Value = -2188.7186  # M132: -2475.799999999999 -> -2188.7186
$ws.Cells.Item(132, 14).Value = -9836.2142  # N132: -12867.8 -> -9836.2142

$ws = $wb.Sheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 970.2195  # H2: 1015.38464 -> 970.2195
$ws.Cells.Item(2, 9).Value = 882.5357  # I2: 943.53845 -> 882.5357
$ws.Cells.Item(2, 11).Value = 882.5357  # K2: 943.53845 -> 882.5357
$ws.Cells.Item(2, 13).Value = -769.5357  # M2: -830.53845 -> -769.5357
$ws.Cells.Item(32, 8).Value = 5623.6387  # H32: 5621.169 -> 5623.6387
$ws.Cells.Item(32, 9).Value = 4927.1714  # I32: 4914.536 -> 4927.1714
$ws.Cells.Item(32, 11).Value = 4927.1714  # K32: 4914.536 -> 4927.1714
$ws.Cells.Item(32, 13).Value = -4640.1714  # M32: -4627.536 -> -4640.1714
$ws.Cells.Item(45, 8).Value = 2166941.2  # H45: 2166983.2 -> 2166941.2
$ws.Cells.Item(45, 9).Value = 3789793.5  # I45: 3789867 -> 3789793.5
$ws.Cells.Item(45, 11).Value = 3789793.5  # K45: 3789867 -> 3789793.5
$ws.Cells.Item(45, 13).Value = -3789416.5  # M45: -3789490 -> -3789416.5
$ws.Cells.Item(49, 8).Value = 30000  # H49: 16525 -> 30000
$ws.Cells.Item(49, 10).Value = 30000  # J49: 16525 -> 30000
$ws.Cells.Item(49, 12).Value = 30000  # L49: 16525 -> 30000
$ws.Cells.Item(49, 14).Value = -30520  # N49: -17045 -> -30520
$ws.Cells.Item(88, 8).Value = 5034.5625  # H88: 5510.5713 -> 5034.5625
$ws.Cells.Item(88, 9).Value = 8758.857  # I88: 13976.5 -> 8758.857
$ws.Cells.Item(88, 10).Value = 2137.889  # J88: 2124.2 -> 2137.889
$ws.Cells.Item(88, 11).Value = 8758.857  # K88: 13976.5 -> 8758.857
$ws.Cells.Item(88, 12).Value = 2137.889  # L88: 2124.2 -> 2137.889
$ws.Cells.Item(88, 13).Value = -8352.857  # M88: -13570.5 -> -8352.857
$ws.Cells.Item(88, 14).Value = -2949.889  # N88: -2936.2 -> -2949.889
$ws.Cells.Item(91, 8).Value = 5034.5625  # H91: 5510.5713 -> 5034.5625
$ws.Cells.Item(91, 9).Value = 8758.857  # I91: 13976.5 -> 8758.857
$ws.Cells.Item(91, 10).Value = 2137.889  # J91: 2124.2 -> 2137.889
$ws.Cells.Item(91, 11).Value = 8758.857  # K91: 13976.5 -> 8758.857
$ws.Cells.Item(91, 12).Value = 2137.889  # L91: 2124.2 -> 2137.889
$ws.Cells.Item(91, 13).Value = -7354.857  # M91: -12572.5 -> -7354.857
$ws.Cells.Item(91, 14).Value = -4945.889  # N91: -4932.2 -> -4945.889
$ws.Cells.Item(97, 8).Value = 1202.1177  # H97: 1401.75 -> 1202.1177
$ws.Cells.Item(97, 9).Value = 736.1111  # I97: 782 -> 736.1111
$ws.Cells.Item(97, 10).Value = 1726.375  # J97: 1844.4286 -> 1726.375
$ws.Cells.Item(97, 11).Value = 736.1111  # K97: 782 -> 736.1111
$ws.Cells.Item(97, 12).Value = 1726.375  # L97: 1844.4286 -> 1726.375
$ws.Cells.Item(97, 13).Value = -240.1111  # M97: -286 -> -240.1111
$ws.Cells.Item(97, 14).Value = -2718.375  # N97: -2836.4286 -> -2718.375
$ws.Cells.Item(102, 8).Value = 1887.8422  # H102: 1942.1666 -> 1887.8422
$ws.Cells.Item(102, 9).Value = 1655.3846  # I102: 1717.5 -> 1655.3846
$ws.Cells.Item(102, 11).Value = 1655.3846  # K102: 1717.5 -> 1655.3846
$ws.Cells.Item(102, 13).Value = -33.38460000000009  # M102: -95.5 -> -33.38460000000009
$ws.Cells.Item(116, 8).Value = 970.2195  # H116: 1015.38464 -> 970.2195
$ws.Cells.Item(116, 9).Value = 882.5357  # I116: 943.53845 -> 882.5357
$ws.Cells.Item(116, 11).Value = 882.5357  # K116: 943.53845 -> 882.5357
$ws.Cells.Item(116, 13).Value = 1411.4643  # M116: 1350.46155 -> 1411.4643
$ws.Cells.Item(132, 8).Value = 7856.825  # H132: 8387.892 -> 7856.825
$ws.Cells.Item(132, 9).Value = 8170.8  # I132: 8203.666999999999 -> 8170.8
$ws.Cells.Item(132, 10).Value = 7668.44  # J132: 8513.5 -> 7668.44
$ws.Cells.Item(132, 11).Value = 24512.4  # K132: 24611.001 -> 24512.4
$ws.Cells.Item(132, 12).Value = 23005.32  # L132: 25540.5 -> 23005.32
$ws.Cells.Item(132, 13).Value = -21982.4  # M132: -22081.001 -> -21982.4
$ws.Cells.Item(132, 14).Value = -28065.32  # N132: -30600.5 -> -28065.32

$ws = $wb.Sheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 970.2195  # H3: 1015.38464 -> 970.2195
$ws.Cells.Item(3, 9).Value = 882.5357  # I3: 943.53845 -> 882.5357
$ws.Cells.Item(3, 11).Value = 882.5357  # K3: 943.53845 -> 882.5357
$ws.Cells.Item(3, 13).Value = -768.5357  # M3: -829.53845 -> -768.5357
$ws.Cells.Item(86, 8).Value = 10103189  # H86: 10103171 -> 10103189
$ws.Cells.Item(86, 9).Value = 10103189  # I86: 10754857 -> 10103189
$ws.Cells.Item(86, 10).Value = 0  # J86: 2050 -> 0
$ws.Cells.Item(86, 11).Value = 10103189  # K86: 10754857 -> 10103189
$ws.Cells.Item(86, 12).Value = 0  # L86: 2050 -> 0
$ws.Cells.Item(86, 13).Value = -10102066  # M86: -10753734 -> -10102066
$ws.Cells.Item(86, 14).ClearContents()  # N86: -4296 -> (removed)
$ws.Cells.Item(89, 8).Value = 10103189  # H89: 10103171 -> 10103189
$ws.Cells.Item(89, 9).Value = 10103189  # I89: 10754857 -> 10103189
$ws.Cells.Item(89, 10).Value = 0  # J89: 2050 -> 0
$ws.Cells.Item(89, 11).Value = 50515945  # K89: 53774285 -> 50515945
$ws.Cells.Item(89, 12).Value = 0  # L89: 10250 -> 0
$ws.Cells.Item(89, 13).Value = -50510329  # M89: -53768669 -> -50510329
$ws.Cells.Item(89, 14).ClearContents()  # N89: -21482 -> (removed)

$ws = $wb.Sheets.Item("CRP")
$ws.Cells.Item(10, 8).Value = 51610  # H10: 51402.2 -> 51610
$ws.Cells.Item(10, 9).Value = 285.33334  # I10: 499.5 -> 285.33334
$ws.Cells.Item(10, 10).Value = 82404.8  # J10: 64127.875 -> 82404.8
$ws.Cells.Item(10, 11).Value = 285.33334  # K10: 499.5 -> 285.33334
$ws.Cells.Item(10, 12).Value = 82404.8  # L10: 64127.875 -> 82404.8
$ws.Cells.Item(10, 13).Value = -146.33334  # M10: -360.5 -> -146.33334
$ws.Cells.Item(10, 14).Value = -82682.8  # N10: -64405.875 -> -82682.8
$ws.Cells.Item(22, 8).Value = 174.44444  # H22: 183.75 -> 174.44444
$ws.Cells.Item(22, 9).Value = 144  # I22: 155 -> 144
$ws.Cells.Item(22, 11).Value = 144  # K22: 155 -> 144
$ws.Cells.Item(22, 13).Value = 206  # M22: 195 -> 206
$ws.Cells.Item(25, 8).Value = 7900  # H25: 0 -> 7900
$ws.Cells.Item(25, 9).Value = 7900  # I25: 0 -> 7900
$ws.Cells.Item(25, 11).Value = 7900  # K25: 0 -> 7900
$ws.Cells.Item(25, 13).Value = -7726  # M25: None -> -7726
$ws.Cells.Item(33, 8).Value = 0  # H33: 1220.25 -> 0
$ws.Cells.Item(33, 9).Value = 0  # I33: 1220.25 -> 0
$ws.Cells.Item(33, 11).Value = 0  # K33: 1220.25 -> 0
$ws.Cells.Item(33, 13).ClearContents()  # M33: -841.25 -> (removed)
$ws.Cells.Item(58, 8).Value = 1685861.1  # H58: 1685926.1 -> 1685861.1
$ws.Cells.Item(58, 9).Value = 2526814.5  # I58: 2599025 -> 2526814.5
$ws.Cells.Item(58, 10).Value = 3954.111  # J58: 3901.7896 -> 3954.111
$ws.Cells.Item(58, 11).Value = 2526814.5  # K58: 2599025 -> 2526814.5
$ws.Cells.Item(58, 12).Value = 3954.111  # L58: 3901.7896 -> 3954.111
$ws.Cells.Item(58, 13).Value = -2526611.5  # M58: -2598822 -> -2526611.5
$ws.Cells.Item(58, 14).Value = -4360.111  # N58: -4307.7896 -> -4360.111
$ws.Cells.Item(136, 8).Value = 1685861.1  # H136: 1685926.1 -> 1685861.1
$ws.Cells.Item(136, 9).Value = 2526814.5  # I136: 2599025 -> 2526814.5
$ws.Cells.Item(136, 10).Value = 3954.111  # J136: 3901.7896 -> 3954.111
$ws.Cells.Item(136, 11).Value = 7580443.5  # K136: 7797075 -> 7580443.5
$ws.Cells.Item(136, 12).Value = 11862.333  # L136: 11705.3688 -> 11862.333
$ws.Cells.Item(136, 13).Value = -7577893.5  # M136: -7794525 -> -7577893.5
$ws.Cells.Item(136, 14).Value = -16962.333  # N136: -16805.3688 -> -16962.333

$ws = $wb.Sheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 150  # H7: 186 -> 150
$ws.Cells.Item(7, 9).Value = 140  # I7: 166.66667 -> 140
$ws.Cells.Item(7, 10).Value = 200  # J7: 215 -> 200
$ws.Cells.Item(7, 11).Value = 420  # K7: 500.00001 -> 420
$ws.Cells.Item(7, 12).Value = 600  # L7: 645 -> 600
$ws.Cells.Item(7, 13).Value = -308  # M7: -388.00001 -> -308
$ws.Cells.Item(7, 14).Value = -824  # N7: -869 -> -824

$ws = $wb.Sheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 7123.1816  # H80: 7885.5 -> 7123.1816
$ws.Cells.Item(80, 9).Value = 8907.857  # I80: 11771 -> 8907.857
$ws.Cells.Item(80, 11).Value = 8907.857  # K80: 11771 -> 8907.857
$ws.Cells.Item(80, 13).Value = -7909.857  # M80: -10773 -> -7909.857
$ws.Cells.Item(83, 8).Value = 7123.1816  # H83: 7885.5 -> 7123.1816
$ws.Cells.Item(83, 9).Value = 8907.857  # I83: 11771 -> 8907.857
$ws.Cells.Item(83, 11).Value = 44539.285  # K83: 58855 -> 44539.285
$ws.Cells.Item(83, 13).Value = -39547.285  # M83: -53863 -> -39547.285
$ws.Cells.Item(102, 8).Value = 3230.4138  # H102: 3272.5715 -> 3230.4138
$ws.Cells.Item(102, 9).Value = 3557.7778  # I102: 3634.111 -> 3557.7778
$ws.Cells.Item(102, 10).Value = 2694.7273  # J102: 2621.8 -> 2694.7273
$ws.Cells.Item(102, 11).Value = 3557.7778  # K102: 3634.111 -> 3557.7778
$ws.Cells.Item(102, 12).Value = 2694.7273  # L102: 2621.8 -> 2694.7273
$ws.Cells.Item(102, 13).Value = -1935.7778  # M102: -2012.111 -> -1935.7778
$ws.Cells.Item(102, 14).Value = -5938.7273  # N102: -5865.8 -> -5938.7273
$ws.Cells.Item(105, 8).Value = 0  # H105: 60000 -> 0
$ws.Cells.Item(105, 10).Value = 0  # J105: 60000 -> 0
$ws.Cells.Item(105, 12).Value = 0  # L105: 60000 -> 0
$ws.Cells.Item(105, 14).ClearContents()  # N105: -66988 -> (removed)
$ws.Cells.Item(122, 8).Value = 9028.571  # H122: 18200 -> 9028.571
$ws.Cells.Item(122, 9).Value = 11740  # I122: 26050 -> 11740
$ws.Cells.Item(122, 10).Value = 2250  # J122: 2500 -> 2250
$ws.Cells.Item(122, 11).Value = 35220  # K122: 78150 -> 35220
$ws.Cells.Item(122, 12).Value = 6750  # L122: 7500 -> 6750
$ws.Cells.Item(122, 13).Value = -32770  # M122: -75700 -> -32770
$ws.Cells.Item(122, 14).Value = -11650  # N122: -12400 -> -11650
$ws.Cells.Item(132, 8).Value = 3091.5652  # H132: 2320.4243 -> 3091.5652
$ws.Cells.Item(132, 9).Value = 2991.8462  # I132: 1880.091 -> 2991.8462
$ws.Cells.Item(132, 10).Value = 3221.2  # J132: 3201.0908 -> 3221.2
$ws.Cells.Item(132, 11).Value = 8975.5386  # K132: 5640.272999999999 -> 8975.5386
$ws.Cells.Item(132, 12).Value = 9663.599999999999  # L132: 9603.2724 -> 9663.599999999999
$ws.Cells.Item(132, 13).Value = -6445.5386  # M132: -3110.272999999999 -> -6445.5386
$ws.Cells.Item(132, 14).Value = -14723.6  # N132: -14663.2724 -> -14723.6

$ws = $wb.Sheets.Item("LTW")
$ws.Cells.Item(41, 8).Value = 0  # H41: 10000 -> 0
$ws.Cells.Item(41, 9).Value = 0  # I41: 10000 -> 0
$ws.Cells.Item(41, 11).Value = 0  # K41: 10000 -> 0
$ws.Cells.Item(41, 13).ClearContents()  # M41: -9562 -> (removed)
$ws.Cells.Item(47, 8).Value = 17995  # H47: 18999 -> 17995
$ws.Cells.Item(47, 10).Value = 17995  # J47: 18999 -> 17995
$ws.Cells.Item(47, 12).Value = 17995  # L47: 18999 -> 17995
$ws.Cells.Item(47, 14).Value = -18975  # N47: -19979 -> -18975
$ws.Cells.Item(52, 8).Value = 17995  # H52: 18999 -> 17995
$ws.Cells.Item(52, 10).Value = 17995  # J52: 18999 -> 17995
$ws.Cells.Item(52, 12).Value = 17995  # L52: 18999 -> 17995
$ws.Cells.Item(52, 14).Value = -18461  # N52: -19465 -> -18461
$ws.Cells.Item(68, 8).Value = 1216.6666  # H68: 1433.3334 -> 1216.6666
$ws.Cells.Item(68, 9).Value = 1300  # I68: 1457.1428 -> 1300
$ws.Cells.Item(68, 10).Value = 966.6667  # J68: 1350 -> 966.6667
$ws.Cells.Item(68, 11).Value = 1300  # K68: 1457.1428 -> 1300
$ws.Cells.Item(68, 12).Value = 966.6667  # L68: 1350 -> 966.6667
$ws.Cells.Item(68, 13).Value = -551  # M68: -708.1428000000001 -> -551
$ws.Cells.Item(68, 14).Value = -2464.6667  # N68: -2848 -> -2464.6667
$ws.Cells.Item(71, 8).Value = 1216.6666  # H71: 1433.3334 -> 1216.6666
$ws.Cells.Item(71, 9).Value = 1300  # I71: 1457.1428 -> 1300
$ws.Cells.Item(71, 10).Value = 966.6667  # J71: 1350 -> 966.6667
$ws.Cells.Item(71, 11).Value = 6500  # K71: 7285.714 -> 6500
$ws.Cells.Item(71, 12).Value = 4833.3335  # L71: 6750 -> 4833.3335
$ws.Cells.Item(71, 13).Value = -2756  # M71: -3541.714 -> -2756
$ws.Cells.Item(71, 14).Value = -12321.3335  # N71: -14238 -> -12321.3335
$ws.Cells.Item(82, 8).Value = 1634  # H82: 1667.6666 -> 1634
$ws.Cells.Item(82, 9).Value = 1332.4445  # I82: 1362.4445 -> 1332.4445
$ws.Cells.Item(82, 10).Value = 2312.5  # J82: 2583.3333 -> 2312.5
$ws.Cells.Item(82, 11).Value = 1332.4445  # K82: 1362.4445 -> 1332.4445
$ws.Cells.Item(82, 12).Value = 2312.5  # L82: 2583.3333 -> 2312.5
$ws.Cells.Item(82, 13).Value = -971.4445000000001  # M82: -1001.4445 -> -971.4445000000001
$ws.Cells.Item(82, 14).Value = -3034.5  # N82: -3305.3333 -> -3034.5
$ws.Cells.Item(85, 8).Value = 1634  # H85: 1667.6666 -> 1634
$ws.Cells.Item(85, 9).Value = 1332.4445  # I85: 1362.4445 -> 1332.4445
$ws.Cells.Item(85, 10).Value = 2312.5  # J85: 2583.3333 -> 2312.5
$ws.Cells.Item(85, 11).Value = 1332.4445  # K85: 1362.4445 -> 1332.4445
$ws.Cells.Item(85, 12).Value = 2312.5  # L85: 2583.3333 -> 2312.5
$ws.Cells.Item(85, 13).Value = -84.44450000000006  # M85: -114.4445000000001 -> -84.44450000000006
$ws.Cells.Item(85, 14).Value = -4808.5  # N85: -5079.3333 -> -4808.5

$ws = $wb.Sheets.Item("WVR")
$ws.Cells.Item(8, 8).Value = 10004  # H8: 11000 -> 10004
$ws.Cells.Item(8, 10).Value = 10004  # J8: 11000 -> 10004
$ws.Cells.Item(8, 12).Value = 10004  # L8: 11000 -> 10004
$ws.Cells.Item(8, 14).Value = -10284  # N8: -11280 -> -10284
$ws.Cells.Item(15, 8).Value = 8631.579  # H15: 8888.888999999999 -> 8631.579
$ws.Cells.Item(64, 8).Value = 40114  # H64: 36114 -> 40114
$ws.Cells.Item(64, 10).Value = 40114  # J64: 36114 -> 40114
$ws.Cells.Item(64, 12).Value = 40114  # L64: 36114 -> 40114
$ws.Cells.Item(64, 14).Value = -40610  # N64: -36610 -> -40610
$ws.Cells.Item(67, 8).Value = 40114  # H67: 36114 -> 40114
$ws.Cells.Item(67, 10).Value = 40114  # J67: 36114 -> 40114
$ws.Cells.Item(67, 12).Value = 40114  # L67: 36114 -> 40114
$ws.Cells.Item(67, 14).Value = -41830  # N67: -37830 -> -41830
$ws.Cells.Item(133, 8).Value = 44610.715  # H133: 53285.715 -> 44610.715
$ws.Cells.Item(133, 10).Value = 44610.715  # J133: 53285.715 -> 44610.715
$ws.Cells.Item(133, 12).Value = 44610.715  # L133: 53285.715 -> 44610.715
$ws.Cells.Item(133, 14).Value = -54730.715  # N133: -63405.715 -> -54730.715
